$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting existing D:K data to E:L
$ws.Columns("D:D").Insert()

# Restore number/date formatting for the new column D by copying formats from column E
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new FY2018 column (D) and revise FY2017/FY2016 figures (E, F) per updated financials
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 12351000
$ws.Range("E8").Value = 12261000
$ws.Range("F8").Value = 13679000
$ws.Range("D9").Value = 5165000
$ws.Range("E9").Value = 6442000
$ws.Range("F9").Value = 5309000
$ws.Range("D10").Value = 7186000
$ws.Range("E10").Value = 5819000
$ws.Range("F10").Value = 8370000
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 49000
$ws.Range("E14").Value = 44000
$ws.Range("D15").Value = 49000
$ws.Range("D17").Value = 9428000
$ws.Range("E17").Value = 11009000
$ws.Range("F17").Value = 9630000
$ws.Range("D18").Value = 2923000
$ws.Range("E18").Value = 1252000
$ws.Range("F18").Value = 4049000
$ws.Range("D20").Value = -231000
$ws.Range("E20").Value = 214000
$ws.Range("F20").Value = -1848000
$ws.Range("D21").Value = 2949000
$ws.Range("E21").Value = 1865000
$ws.Range("F21").Value = 2959000
$ws.Range("D22").Value = 231000
$ws.Range("D23").Value = 2461000
$ws.Range("E23").Value = 1306000
$ws.Range("F23").Value = 2027000
$ws.Range("D24").Value = 380000
$ws.Range("E24").Value = 125000
$ws.Range("F24").Value = 378000
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 2081000
$ws.Range("E26").Value = 1181000
$ws.Range("F26").Value = 1649000
$ws.Range("D27").Value = 1747000
$ws.Range("E27").Value = 758000
$ws.Range("F27").Value = 1254000
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 73000
$ws.Range("E29").Value = 76000
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 231000
$ws.Range("E32").Value = -214000
$ws.Range("F32").Value = 1848000
$ws.Range("D33").Value = 1820000
$ws.Range("E33").Value = 834000
$ws.Range("F33").Value = 1254000
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 1820000
$ws.Range("E35").Value = 834000
$ws.Range("F35").Value = 1254000
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 4469000
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 7104000
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 81281000
$ws.Range("E47").Value = 82706000
$ws.Range("D48").Value = 0
$ws.Range("D49").Value = 4961000
$ws.Range("E49").Value = 9810000
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 7932000
$ws.Range("E52").Value = 7599000
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 220797000
$ws.Range("E54").Value = 235615000
$ws.Range("D57").Value = 0
$ws.Range("D58").Value = 2609000
$ws.Range("E58").Value = 5336000
$ws.Range("D59").Value = 5032000
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 10899000
$ws.Range("E61").Value = 9673000
$ws.Range("D62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 206969000
$ws.Range("E66").Value = 222246000
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 13989000
$ws.Range("E72").Value = 12225000
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 13828000
$ws.Range("E76").Value = 13369000
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 1820000
$ws.Range("E81").Value = 834000
$ws.Range("F81").Value = 1254000
$ws.Range("D83").Value = 257000
$ws.Range("E83").Value = 399000
$ws.Range("F83").Value = 758000
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 61000
$ws.Range("E89").Value = -243000
$ws.Range("D91").Value = -123000
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -2049000
$ws.Range("D96").Value = -157000
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 1655000
$ws.Range("E100").Value = 9070000
$ws.Range("D101").Value = -12000
$ws.Range("D102").Value = -345000
